$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row (row 1) with 4 new snapshot columns: DW:DZ
# Copy style from the previous snapshot block (DS1:DV1) so the new headers
# keep the bold/centered/bordered header formatting.
$ws.Range("DS1:DV1").Copy($ws.Range("DW1:DZ1"))
$ws.Range("DW1").Value = "Daily as on Jun 7, 16:00"
$ws.Range("DX1").Value = "Weekly as on Jun 7, 16:00"
$ws.Range("DY1").Value = "Monthly as on Jun 7, 16:00"
$ws.Range("DZ1").Value = "Closing as on Jun 7, 16:00"

# Fill in the new Jun 7, 16:00 snapshot data for every stock row (2-39)
$ws.Range("DW2").Value = 75.48999999999999
$ws.Range("DX2").Value = 63.17
$ws.Range("DY2").Value = 66.31
$ws.Range("DZ2").Value = 2227.4
$ws.Range("DW3").Value = 69.77
$ws.Range("DX3").Value = 72.23
$ws.Range("DY3").Value = 67.92
$ws.Range("DZ3").Value = 345.5
$ws.Range("DW4").Value = 56.05
$ws.Range("DX4").Value = 60.09
$ws.Range("DY4").Value = 67.91
$ws.Range("DZ4").Value = 5729.75
$ws.Range("DW5").Value = 61.97
$ws.Range("DX5").Value = 69.64
$ws.Range("DY5").Value = 70.29000000000001
$ws.Range("DZ5").Value = 11828.55
$ws.Range("DW6").Value = 55.77
$ws.Range("DX6").Value = 58.89
$ws.Range("DY6").Value = 63.12
$ws.Range("DZ6").Value = 1499.85
$ws.Range("DW7").Value = 55.54
$ws.Range("DX7").Value = 53.56
$ws.Range("DY7").Value = 62.09
$ws.Range("DZ7").Value = 680.1
$ws.Range("DW8").Value = 67.31
$ws.Range("DX8").Value = 57.9
$ws.Range("DY8").Value = 61.03
$ws.Range("DZ8").Value = 3032.7
$ws.Range("DW9").Value = 71.65000000000001
$ws.Range("DX9").Value = 55.77
$ws.Range("DY9").Value = 54.29
$ws.Range("DZ9").Value = 7275.65
$ws.Range("DW10").Value = 68.41
$ws.Range("DX10").Value = 70.08
$ws.Range("DY10").Value = 68.19
$ws.Range("DZ10").Value = 4249.75
$ws.Range("DW11").Value = 60.12
$ws.Range("DX11").Value = 63.74
$ws.Range("DY11").Value = 63.82
$ws.Range("DZ11").Value = 60.4
$ws.Range("DW12").Value = 76.65000000000001
$ws.Range("DX12").Value = 72.87
$ws.Range("DY12").Value = 58.82
$ws.Range("DZ12").Value = 414.7
$ws.Range("DW13").Value = 67.09999999999999
$ws.Range("DX13").Value = 71.63
$ws.Range("DY13").Value = 71.81
$ws.Range("DZ13").Value = 242.5
$ws.Range("DW14").Value = 60.41
$ws.Range("DX14").Value = 71.73999999999999
$ws.Range("DY14").Value = 65.41
$ws.Range("DZ14").Value = 2369.8
$ws.Range("DW15").Value = 50.38
$ws.Range("DX15").Value = 50.76
$ws.Range("DY15").Value = 59.04
$ws.Range("DZ15").Value = 536.45
$ws.Range("DW16").Value = 68.76000000000001
$ws.Range("DX16").Value = 63.3
$ws.Range("DY16").Value = 46.92
$ws.Range("DZ16").Value = 156.25
$ws.Range("DW17").Value = 72.56
$ws.Range("DX17").Value = 77.23
$ws.Range("DY17").Value = 84.06999999999999
$ws.Range("DZ17").Value = 548.25
$ws.Range("DW18").Value = 70.86
$ws.Range("DX18").Value = 60.6
$ws.Range("DY18").Value = 44
$ws.Range("DZ18").Value = 247.8
$ws.Range("DW19").Value = 73.81
$ws.Range("DX19").Value = 53.53
$ws.Range("DY19").Value = 43.28
$ws.Range("DZ19").Value = 10.1
$ws.Range("DW20").Value = 61.46
$ws.Range("DX20").Value = 59.16
$ws.Range("DY20").Value = 50.76
$ws.Range("DZ20").Value = 1025.2
$ws.Range("DW21").Value = 75.15000000000001
$ws.Range("DX21").Value = 69.26000000000001
$ws.Range("DY21").Value = 53.1
$ws.Range("DZ21").Value = 114.7
$ws.Range("DW22").Value = 78.54000000000001
$ws.Range("DX22").Value = 69.56999999999999
$ws.Range("DY22").Value = 72.04000000000001
$ws.Range("DZ22").Value = 2094.8
$ws.Range("DW23").Value = 53.6
$ws.Range("DX23").Value = 52.26
$ws.Range("DY23").Value = 47.88
$ws.Range("DZ23").Value = 211.45
$ws.Range("DW24").Value = 75.73999999999999
$ws.Range("DX24").Value = 68.43000000000001
$ws.Range("DY24").Value = 63.55
$ws.Range("DZ24").Value = 1557.8
$ws.Range("DW25").Value = 58.79
$ws.Range("DX25").Value = 49.25
$ws.Range("DY25").Value = 48.85
$ws.Range("DZ25").Value = 170.7
$ws.Range("DW26").Value = 62.07
$ws.Range("DX26").Value = 54.3
$ws.Range("DY26").Value = 54.93
$ws.Range("DZ26").Value = 165.4
$ws.Range("DW27").Value = 81.91
$ws.Range("DX27").Value = 69.84
$ws.Range("DY27").Value = 69.62
$ws.Range("DZ27").Value = 1513.15
$ws.Range("DW28").Value = 69.53
$ws.Range("DX28").Value = 60.75
$ws.Range("DY28").Value = 58
$ws.Range("DZ28").Value = 125.75
$ws.Range("DW29").Value = 59.92
$ws.Range("DX29").Value = 64.95999999999999
$ws.Range("DY29").Value = 61.91
$ws.Range("DZ29").Value = 988.5
$ws.Range("DW30").Value = 73.61
$ws.Range("DX30").Value = 73.68000000000001
$ws.Range("DY30").Value = 67.2
$ws.Range("DZ30").Value = 432.25
$ws.Range("DW31").Value = 66.68000000000001
$ws.Range("DX31").Value = 71.31999999999999
$ws.Range("DY31").Value = 72.64
$ws.Range("DZ31").Value = 113.35
$ws.Range("DW32").Value = 59.03
$ws.Range("DX32").Value = 58.85
$ws.Range("DY32").Value = 69.43000000000001
$ws.Range("DZ32").Value = 3183.2
$ws.Range("DW33").Value = 60.86
$ws.Range("DX33").Value = 62.01
$ws.Range("DY33").Value = 61.87
$ws.Range("DZ33").Value = 748.2
$ws.Range("DW34").Value = 67.22
$ws.Range("DX34").Value = 67.52
$ws.Range("DY34").Value = 52.97
$ws.Range("DZ34").Value = 10.95
$ws.Range("DW35").Value = 66.81
$ws.Range("DX35").Value = 69.88
$ws.Range("DY35").Value = 72.95999999999999
$ws.Range("DZ35").Value = 2933.1
$ws.Range("DW36").Value = 75.17
$ws.Range("DX36").Value = 73.83
$ws.Range("DY36").Value = 83.22
$ws.Range("DZ36").Value = 878.6
$ws.Range("DW37").Value = 50.59
$ws.Range("DX37").Value = 63.89
$ws.Range("DY37").Value = 66.95
$ws.Range("DZ37").Value = 675.35
$ws.Range("DW38").Value = 58.48
$ws.Range("DX38").Value = 78.09999999999999
$ws.Range("DY38").Value = 79.64
$ws.Range("DZ38").Value = 1128.7
$ws.Range("DW39").Value = 50.3
$ws.Range("DX39").Value = 68.63
$ws.Range("DY39").Value = 84.43000000000001
$ws.Range("DZ39").Value = 1754.15
